$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the cells can be updated.
$ws.Unprotect("lido")

$newNote = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."
$ws.Range("A11").Value = $newNote

$ws.Range("D2").Value = 0.4948073089411785
$ws.Range("E2").Value = -0.005629838142153232

$ws.Range("D3").Value = 0.241843363778022
$ws.Range("E3").Value = -0.02337742233159013

$ws.Range("D4").Value = 0.09909164508156437
$ws.Range("E4").Value = -0.02440251572327046

$ws.Range("D5").Value = 0.103780617804801
$ws.Range("E5").Value = -0.01401201029453814

$ws.Range("D6").Value = 0.03141018183717512
$ws.Range("E6").Value = -0.01795162509448234

$ws.Range("D7").Value = 0.02906688255725889
$ws.Range("E7").Value = -0.02687930518119197

$ws.Range("D8").Value = 0.9999999999999999
$ws.Range("E8").Value = -0.01365678144163962

# Restore sheet protection to its prior (protected) state.
$ws.Protect("lido")
